$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sheet3" (the active/intern-hours sheet) gains a third column that records
# whether a time-off entry is still enabled, and the stray inline-string
# helper row for the intern's name is replaced by a normal table row.

# New header + values for column C ("Enabled").
$ws.Range("C1").Value = "Enabled"
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C1:C3").NumberFormat = "0"
$ws.Columns("C").ColumnWidth = 8.43

# Row 3's name cell used to be a one-off inline string outside the table;
# re-enter the same text through the Name column so it becomes a normal
# (shared-string) cell like A2, completing the Name/Date/Enabled row for
# the intern's second day.
$ws.Range("A3").Value = "Austin Danaj"

# Move the selection to match where Excel left the cursor after the edit.
$ws.Range("B4").Select()
